$d = $word.ActiveDocument

# Remove the trailing stand-alone space run in the first paragraph
# (it sits right before the paragraph mark).
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$spaceRange = $d.Range($r1.End - 2, $r1.End - 1)
$spaceRange.Delete()

# Update the placeholder id text in the (now single) run.
$d.Content.Find.Execute("**ID__AFFARS_5325_topic_17__ID**", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_SUBPART_5325_73__ID**", 2)

# Widen the left indent of the first paragraph (120 -> 225 twips = 6pt -> 11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Add a (line-less) paragraph border with 5-twip spacing on every side.
$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5
